# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the existing "总计" sheet,
#    and populate it with the per-fund holdings detail for 2022-Q1.
# 2. Update the "总计" (grand total) summary sheet so it gets a new first
#    data row for "2022-Q1" (existing rows shift down by one).

function Set-TextCell($ws, $row, $col, $val) {
    # Force the cell to be stored as text even when the value looks like a
    # number (fund codes, percentages recorded as plain text, etc.) - mirrors
    # how the source workbook stores these columns as inlineStr/text.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

function Set-NumberCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q1" sheet right before "总计"
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($beforeSheet)
$q1.Name = "2022-Q1"

# NOTE: the sheet reference fetched before Add() can become aliased to the
# freshly-inserted sheet once the worksheets collection is mutated, so the
# "总计" sheet must be re-resolved by name AFTER the Add() call.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row
Set-TextCell $q1 1 2 "基金代码"
Set-TextCell $q1 1 3 "基金名称"
Set-TextCell $q1 1 4 "基金规模"
Set-TextCell $q1 1 5 "股票总仓位"
Set-TextCell $q1 1 6 "仓位占比"
Set-TextCell $q1 1 7 "持有市值(亿元)"
Set-TextCell $q1 1 8 "仓位排名"

# Data rows: (row index col A, code, name, scale, position, ratio, marketValue, rank)
$q1Rows = @(
    @(0, "010845", "泰达宏利波控回报12个月持有期混合", "13.29", "23.79", "0.64", "0.0851", 10),
    @(1, "393001", "中海优势精选灵活配置混合",          "1.58",  "78.65", "5.03", "0.0795", 10),
    @(2, "163110", "申万菱信量化小盘股票(LOF)",          "5.68",  "92.25", "1.31", "0.0744", 5),
    @(3, "162205", "泰达宏利风险预算混合",                "1.34",  "27.82", "0.76", "0.0102", 9)
)

$r = 2
foreach ($row in $q1Rows) {
    Set-NumberCell $q1 $r 1 $row[0]
    Set-TextCell   $q1 $r 2 $row[1]
    Set-TextCell   $q1 $r 3 $row[2]
    Set-TextCell   $q1 $r 4 $row[3]
    Set-TextCell   $q1 $r 5 $row[4]
    Set-TextCell   $q1 $r 6 $row[5]
    Set-TextCell   $q1 $r 7 $row[6]
    Set-NumberCell $q1 $r 8 $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet - add a new "2022-Q1" row at the top of
# the data (row 2), pushing the previous rows down by one.
# ---------------------------------------------------------------------
Set-TextCell $totalSheet 1 2 "日期"
Set-TextCell $totalSheet 1 3 "持有数量(只)"
Set-TextCell $totalSheet 1 4 "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 4,  0.25),
    @(1, "2021-Q4", 1,  0.08),
    @(2, "2021-Q3", 12, 2.27),
    @(3, "2021-Q2", 6,  0.84),
    @(4, "2021-Q1", 3,  0.34)
)

$r = 2
foreach ($row in $totalRows) {
    Set-NumberCell $totalSheet $r 1 $row[0]
    Set-TextCell   $totalSheet $r 2 $row[1]
    Set-NumberCell $totalSheet $r 3 $row[2]
    Set-NumberCell $totalSheet $r 4 $row[3]
    $r = $r + 1
}
